$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.909168386128414
$ws.Range("C2").Value = 0.03105757396912168
$ws.Range("E2").Value = 0.07852987395450661
$ws.Range("F2").Value = 3.742124059222519
$ws.Range("G2").Value = 0.002618631557067928
$ws.Range("J2").Value = 0.2199797100253313
$ws.Range("K2").Value = 1.363528138879758
$ws.Range("L2").Value = 0.2877771645403016
$ws.Range("N2").Value = 3.858137254446149
$ws.Range("B3").Value = 1.868268787627954
$ws.Range("C3").Value = 0.02705245307993209
$ws.Range("E3").Value = 0.07847626374588224
$ws.Range("F3").Value = 3.72276296164921
$ws.Range("G3").Value = 0.002623036941766184
$ws.Range("J3").Value = 0.2207590449109098
$ws.Range("K3").Value = 1.322264614845722
$ws.Range("L3").Value = 0.2839071752102456
$ws.Range("N3").Value = 3.865879027646415
$ws.Range("B4").Value = 1.84422793003165
$ws.Range("C4").Value = 0.02458985849845874
$ws.Range("E4").Value = 0.07847429984379239
$ws.Range("F4").Value = 3.712388858739118
$ws.Range("G4").Value = 0.002625885012283175
$ws.Range("J4").Value = 0.2212656576236984
$ws.Range("K4").Value = 1.297738820828926
$ws.Range("L4").Value = 0.2816884988424917
$ws.Range("N4").Value = 3.871439607836891
$ws.Range("B5").Value = 1.834700585638245
$ws.Range("C5").Value = 0.02358537901130831
$ws.Range("E5").Value = 0.07848131024403671
$ws.Range("F5").Value = 3.708541910237088
$ws.Range("G5").Value = 0.002627081736213858
$ws.Range("J5").Value = 0.2214791720152061
$ws.Range("K5").Value = 1.287947894092923
$ws.Range("L5").Value = 0.2808240558563355
$ws.Range("N5").Value = 3.873908350958061
$ws.Range("B6").Value = 1.833134855451192
$ws.Range("C6").Value = 0.02341852586688731
$ws.Range("E6").Value = 0.0784829468698014
$ws.Range("F6").Value = 3.707926112482255
$ws.Range("G6").Value = 0.002627282635635386
$ws.Range("J6").Value = 0.2215150526703979
$ws.Range("K6").Value = 1.286334405712637
$ws.Range("L6").Value = 0.2806829149873664
$ws.Range("N6").Value = 3.874330524669034
$ws.Range("B7").Value = 1.844098349705916
$ws.Range("C7").Value = 0.02457631570862873
$ws.Range("E7").Value = 0.07847436272648523
$ws.Range("F7").Value = 3.712335436554923
$ws.Range("G7").Value = 0.002625901005481355
$ws.Range("J7").Value = 0.2212685085454336
$ws.Range("K7").Value = 1.297605953118762
$ws.Range("L7").Value = 0.2816766798713317
$ws.Range("N7").Value = 3.871472081470941
$ws.Range("B8").Value = 1.894843773846816
$ws.Range("C8").Value = 0.02967726648671487
$ws.Range("E8").Value = 0.07850498101485393
$ws.Range("F8").Value = 3.735134142287379
$ws.Range("G8").Value = 0.002620120894073119
$ws.Range("J8").Value = 0.2202425930644338
$ws.Range("K8").Value = 1.349132189209115
$ws.Range("L8").Value = 0.2864101447966192
$ws.Range("N8").Value = 3.860638986865524
$ws.Range("B9").Value = 2.002865838941091
$ws.Range("C9").Value = 0.03965741933733113
$ws.Range("E9").Value = 0.07880950553491495
$ws.Range("F9").Value = 3.791861458883417
$ws.Range("G9").Value = 0.002609916582561193
$ws.Range("J9").Value = 0.2184536995451793
$ws.Range("K9").Value = 1.456619495546363
$ws.Range("L9").Value = 0.2969397351569967
$ws.Range("N9").Value = 3.845808985730997
$ws.Range("B10").Value = 2.087439931090387
$ws.Range("C10").Value = 0.04698275735474056
$ws.Range("E10").Value = 0.07918090734951555
$ws.Range("F10").Value = 3.840887715410361
$ws.Range("G10").Value = 0.002603101161769095
$ws.Range("J10").Value = 0.2172752505033095
$ws.Range("K10").Value = 1.539552273711934
$ws.Range("L10").Value = 0.3054341633006743
$ws.Range("N10").Value = 3.838838336938664
$ws.Range("B11").Value = 2.127051900416802
$ws.Range("C11").Value = 0.05031521922722959
$ws.Range("E11").Value = 0.07938165191775504
$ws.Range("F11").Value = 3.864792115123549
$ws.Range("G11").Value = 0.002600147078683519
$ws.Range("J11").Value = 0.2167686361664813
$ws.Range("K11").Value = 1.578148747993367
$ws.Range("L11").Value = 0.309462837719451
$ws.Range("N11").Value = 3.836522900695272
$ws.Range("B12").Value = 2.142215912978259
$ws.Range("C12").Value = 0.05157726315817968
$ws.Range("E12").Value = 0.07946221519727104
$ws.Range("F12").Value = 3.874074742701936
$ws.Range("G12").Value = 0.002599049358979161
$ws.Range("J12").Value = 0.2165810333071896
$ws.Range("K12").Value = 1.592889759469188
$ws.Range("L12").Value = 0.3110119966566884
$ws.Range("N12").Value = 3.835769400000302
$ws.Range("B13").Value = 2.138942784526307
$ws.Range("C13").Value = 0.0513054528832555
$ws.Range("E13").Value = 0.07944466269530182
$ws.Range("F13").Value = 3.87206530619352
$ws.Range("G13").Value = 0.00259928484349372
$ws.Range("J13").Value = 0.2166212482606298
$ws.Range("K13").Value = 1.589709438910461
$ws.Range("L13").Value = 0.3106773096085078
$ws.Range("N13").Value = 3.835926191356378
$ws.Range("B14").Value = 2.12829616927371
$ws.Range("C14").Value = 0.05041904553411314
$ws.Range("E14").Value = 0.07938818895328303
$ws.Range("F14").Value = 3.865551180913201
$ws.Range("G14").Value = 0.002600056349943232
$ws.Range("J14").Value = 0.2167531169489969
$ws.Range("K14").Value = 1.579358985821273
$ws.Range("L14").Value = 0.3095898156277883
$ws.Range("N14").Value = 3.836458436919301
$ws.Range("B15").Value = 2.12179614618907
$ws.Range("C15").Value = 0.04987611296243699
$ws.Range("E15").Value = 0.07935418836847674
$ws.Range("F15").Value = 3.861591120692339
$ws.Range("G15").Value = 0.002600531642217518
$ws.Range("J15").Value = 0.2168344427405984
$ws.Range("K15").Value = 1.57303536870117
$ws.Range("L15").Value = 0.3089267638812885
$ws.Range("N15").Value = 3.836800518949417
$ws.Range("B16").Value = 2.084874116241167
$ws.Range("C16").Value = 0.04676498377988025
$ws.Range("E16").Value = 0.07916842573034799
$ws.Range("F16").Value = 3.839357759432829
$ws.Range("G16").Value = 0.0026032971526246
$ws.Range("J16").Value = 0.2173089522088354
$ws.Range("K16").Value = 1.5370474314403
$ws.Range("L16").Value = 0.3051741844428193
$ws.Range("N16").Value = 3.839006895598132
$ws.Range("B17").Value = 2.062515359706765
$ws.Range("C17").Value = 0.04485650674784836
$ws.Range("E17").Value = 0.07906259050958653
$ws.Range("F17").Value = 3.82612876834267
$ws.Range("G17").Value = 0.002605031093661836
$ws.Range("J17").Value = 0.2176075982893089
$ws.Range("K17").Value = 1.515192997419632
$ws.Range("L17").Value = 0.3029141850117583
$ws.Range("N17").Value = 3.840579758630227
$ws.Range("B18").Value = 2.04976237829834
$ws.Range("C18").Value = 0.04375880767493356
$ws.Range("E18").Value = 0.07900471084549743
$ws.Range("F18").Value = 3.81867059009673
$ws.Range("G18").Value = 0.002606042186324986
$ws.Range("J18").Value = 0.2177821450598429
$ws.Range("K18").Value = 1.502704771028476
$ws.Range("L18").Value = 0.3016297822892255
$ws.Range("N18").Value = 3.841564939061129
$ws.Range("B19").Value = 2.045462845122358
$ws.Range("C19").Value = 0.04338714387189668
$ws.Range("E19").Value = 0.07898562879584858
$ws.Range("F19").Value = 3.816171271984899
$ws.Range("G19").Value = 0.002606386894154822
$ws.Range("J19").Value = 0.2178417198363256
$ws.Range("K19").Value = 1.498490524689259
$ws.Range("L19").Value = 0.3011975680261685
$ws.Range("N19").Value = 3.841912323693833
$ws.Range("B20").Value = 2.064884393969407
$ws.Range("C20").Value = 0.04505966603618106
$ws.Range("E20").Value = 0.0790735471750672
$ws.Range("F20").Value = 3.82752140993972
$ws.Range("G20").Value = 0.002604845087595987
$ws.Range("J20").Value = 0.2175755198402811
$ws.Range("K20").Value = 1.517510964502691
$ws.Range("L20").Value = 0.3031531632481119
$ws.Range("N20").Value = 3.840403990118716
$ws.Range("B21").Value = 2.131418891685257
$ws.Range("C21").Value = 0.05067940103734259
$ws.Range("E21").Value = 0.07940465349818382
$ws.Range("F21").Value = 3.867458278955212
$ws.Range("G21").Value = 0.002599829172666369
$ws.Range("J21").Value = 0.2167142687901702
$ws.Range("K21").Value = 1.58239576025548
$ws.Range("L21").Value = 0.3099085993681143
$ws.Range("N21").Value = 3.836298754797497
$ws.Range("B22").Value = 2.175857934056125
$ws.Range("C22").Value = 0.05435292539648628
$ws.Range("E22").Value = 0.07964753435011929
$ws.Range("F22").Value = 3.89490323943852
$ws.Range("G22").Value = 0.002596672913044041
$ws.Range("J22").Value = 0.2161761130449307
$ws.Range("K22").Value = 1.625532575630331
$ws.Range("L22").Value = 0.3144611210232569
$ws.Range("N22").Value = 3.834334593282009
$ws.Range("B23").Value = 2.152052583739589
$ws.Range("C23").Value = 0.05239220062453853
$ws.Range("E23").Value = 0.07951548953410992
$ws.Range("F23").Value = 3.880132321716502
$ws.Range("G23").Value = 0.002598346347055024
$ws.Range("J23").Value = 0.216461073628107
$ws.Range("K23").Value = 1.602442681429096
$ws.Range("L23").Value = 0.3120188003091755
$ws.Range("N23").Value = 3.835317035551114
$ws.Range("B24").Value = 2.063813037844056
$ws.Range("C24").Value = 0.04496781925648463
$ws.Range("E24").Value = 0.07906858442505005
$ws.Range("F24").Value = 3.82689133771828
$ws.Range("G24").Value = 0.002604929136515921
$ws.Range("J24").Value = 0.217590013625852
$ws.Range("K24").Value = 1.516462774319479
$ws.Range("L24").Value = 0.3030450747690878
$ws.Range("N24").Value = 3.840483203009896
$ws.Range("B25").Value = 1.972729661735627
$ws.Range("C25").Value = 0.03695926195052834
$ws.Range("E25").Value = 0.07870109306278827
$ws.Range("F25").Value = 3.775226295723385
$ws.Range("G25").Value = 0.002612556870012276
$ws.Range("J25").Value = 0.2189137832545782
$ws.Range("K25").Value = 1.426848013453366
$ws.Range("L25").Value = 0.2939578869588644
$ws.Range("N25").Value = 3.849132733110025
